# Update the "Metadata" sheet (sheet1):
#  - Version bumped 5.0.0 -> 6.0.0
#  - Date bumped to the new publish timestamp
#  - Publisher value filled in ("Alvearie Team")
#  - The old duplicated "Contact" row is turned into the new "Jurisdiction" row,
#    and the extra duplicate "Contact" row is removed entirely.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

$ws.Range("B3").Value = "6.0.0"
$ws.Range("B8").Value = "2022-01-21T20:46:54+00:00"
$ws.Range("B9").Value = "Alvearie Team"
$ws.Range("A10").Value = "Jurisdiction"
$ws.Range("B10").Value = "United States of America"

# Row 11 was a redundant duplicate of the old "Contact" row - remove it so
# everything below shifts up by one (dimension becomes A1:B20).
$ws.Rows.Item(11).Delete()

# Update the "Elements" sheet (sheet2): the root Extension row's Short /
# Definition columns (K2 / L2) now describe this specific extension instead
# of the generic placeholder text.
$ws2 = $wb.Worksheets.Item("Elements")
$ws2.Range("K2").Value = "Employee Cost Center"
$ws2.Range("L2").Value = "Code for the cost center of the employee"
